$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.847.61'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.622.00'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'210.38"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'23.26"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').Value = "'0.0611"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('D11').Value = "'0.0878"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').Value = '1.849.97'
$ws.Range('E12').Value = '  -1.19%  '
$ws.Range('D13').Value = '1.618.21'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').Value = "'4.02"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.89%  '
$ws.Range('D15').Value = "'0.559"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.85%  '
$ws.Range('D16').Value = "'65.17"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.98%  '
$ws.Range('D17').Value = '27.828.75'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = "'228.54"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.14%  '
$ws.Range('D19').Value = '0.0₃0720'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').Value = "'7.59"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').Value = "'1.00"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = "'4.30"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('D23').Value = "'10.07"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.94%  '
$ws.Range('D24').Value = "'2.03"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.79%  '
$ws.Range('D25').Value = "'153.97"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('E26').Value = '  -0.99%  '
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'15.46"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').Value = "'1.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('D31').Value = "'0.0480"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('D32').Value = "'3.41"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.54%  '
$ws.Range('D33').Value = "'3.07"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = '1.384.58'
$ws.Range('E34').Value = '  -2.12%  '
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('E36').Value = '  +11.25%  '
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('D38').Value = "'0.0170"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('D40').Value = "'0.845"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.30%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = "'0.995"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.92%  '
$ws.Range('D43').Value = "'1.85"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('D44').Value = "'5.47"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('D45').Value = "'65.57"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.90%  '
$ws.Range('D46').Value = '1.759.97'
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('D47').Value = "'2.15"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.59%  '
$ws.Range('D48').Value = "'87.60"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.15%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = "'0.101"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.0502"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'7.60"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.12%  '
